$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended to the report: date serial, nuovi pos.,
# somma mobile 7gg., rapporto %
$newRows = @(
    @{ Row = 227; Date = 44301; B = 1; C = 9;  D = 109.3825960136121 },
    @{ Row = 228; Date = 44302; B = 3; C = 12; D = 145.8434613514827 },
    @{ Row = 229; Date = 44303; B = 2; C = 12; D = 145.8434613514827 }
)

foreach ($r in $newRows) {
    $prevRow = $r.Row - 1

    # Match the date-cell formatting of the row above (style index carries
    # the date number format, bold font, border, alignment).
    $ws.Range("A$prevRow").Copy()
    $ws.Range("A$($r.Row)").PasteSpecial(-4122)

    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

$excel.CutCopyMode = $false
